# IETF100-NVO3-IOAM.pptx minor update
#
# Slide 4 ("Section 4: Discussion of the encapsulation approach"),
# shape "Text Placeholder 5":
#   - paragraph "Geneve / tunnel option limits length to 128 bytes...":
#     merge the lone-space run and the following run into a single run
#     that reads " tunnel option limits length to 128 bytes, which
#     limits the range of deployment cases."
#   - paragraph "Hardware-friendly implementation discussion:" becomes
#     "Metadata approach discussion:"
#   - paragraph "Use of the Geneve tunnel option requires ..." keeps its
#     wording, but the leading "Use of the " run is now split into two
#     runs, "Use " and "of the "

$p  = $ppt.ActivePresentation
$s4 = $p.Slides.Item(4)
$sh = $s4.Shapes.Item("Text Placeholder 5")
$tr = $sh.TextFrame.TextRange

# --- "Geneve" + " " + "tunnel option limits ..." -> "Geneve" + " tunnel option limits ..."
$para4 = $tr.Paragraphs(4, 1)
$len4  = $para4.Length
$rest4 = $para4.Characters(7, $len4 - 7)
$rest4.Text = " tunnel option limits length to 128 bytes, which limits the range of deployment cases."

# --- "Hardware-friendly " -> "Metadata approach "; "implementation discussion:" -> "discussion:"
$para5  = $tr.Paragraphs(5, 1)
$run1_5 = $para5.Characters(1, 18)
$run1_5.Text = "Metadata approach "
$para5b = $tr.Paragraphs(5, 1)
$run2_5 = $para5b.Characters(19, 26)
$run2_5.Text = "discussion:"

# --- "Use of the " -> "Use " + "of the " (two runs, same formatting)
$para6  = $tr.Paragraphs(6, 1)
$run1_6 = $para6.Characters(1, 4)
$run1_6.Text = "Use "
$para6b = $tr.Paragraphs(6, 1)
$run2_6 = $para6b.Characters(5, 7)
$run2_6.Text = "of the "
